# save 05/02 / création models et objects
#
# Clean up the "User" / "CartRelation" update-function rows that no longer
# apply, and document the new `password` attribute on the User object.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Style" function block (O:T, row 7): drop the stray updateType() row that
# duplicated the one already documented at row 16.
$ws.Range("P7:R7").ClearContents()

# "User" attribute block (A:F, row 12): document the new `password` field.
$ws.Range("B12").Value = "password"
$ws.Range("C12").Value = "string"
$ws.Range("D12").Value = "mot de passe"

# "User" function block (A:F, row 14): remove the obsolete updateUser() row.
$ws.Range("A14:D14").ClearContents()

# "CartRelation" function block (O:T, row 20): remove the obsolete
# updateCartRelation() row.
$ws.Range("P20:R20").ClearContents()

# Leave the selection where the author ended up editing.
$ws.Range("B12").Select()
